# Updated symbol list on Tue Dec 27 20:29:38 UTC 2022 with GitHub Actions
#
# Applies the per-cell text updates described by the diff. Price values in
# column D are stored as plain text in the workbook (inline strings), so we
# force the cell to Text format before assigning the literal string and then
# restore the cell's style to Normal, which keeps the value as text without
# leaving any visible number-formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# --- Column D (Price) updates -------------------------------------------
Set-TextValue "D2"  "245.92"
Set-TextValue "D3"  "23.97"
Set-TextValue "D4"  "5.364"
Set-TextValue "D5"  "0.05838"
Set-TextValue "D6"  "6.472"
Set-TextValue "D7"  "3.354"
Set-TextValue "D8"  "0.8098"
Set-TextValue "D9"  "0.9192"
Set-TextValue "D11" "0.07404"
Set-TextValue "D12" "0.03109"
Set-TextValue "D14" "0.09373"
Set-TextValue "D15" "3.867"
Set-TextValue "D16" "0.001560"
Set-TextValue "D17" "0.04697"
Set-TextValue "D18" "0.0005987"
Set-TextValue "D19" "0.005962"
Set-TextValue "D20" "0.001245"
Set-TextValue "D21" "0.004691"
Set-TextValue "D22" "0.00008797"
Set-TextValue "D23" "3.595"
Set-TextValue "D25" "0.3183"
Set-TextValue "D26" "0.1319"

# --- Rows 41-43: coin rows shifted up (KickToken/BKEXToken/CEJI reorder) -
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006358"
Set-TextValue "E41" "40KickTokenKICK"

Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1066"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003199"
Set-TextValue "E43" "42CEJICEJI"

# --- Remaining column D (Price) updates ----------------------------------
Set-TextValue "D44" "0.008522"
Set-TextValue "D45" "0.00005246"
Set-TextValue "D47" "0.6858"
Set-TextValue "D48" "0.001837"

Write-Host "Done applying crypto price/symbol updates."
